$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "310.21"
$ws.Range("E2").Value = "-1.18%"
$ws.Range("D3").Value = "37.48"
$ws.Range("E3").Value = "-4.45%"
$ws.Range("D4").Value = "5.078"
$ws.Range("E4").Value = "-0.97%"
$ws.Range("D5").Value = "0.07753"
$ws.Range("E5").Value = "-5.00%"
$ws.Range("D6").Value = "4.346"
$ws.Range("E6").Value = "-2.85%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.891"
$ws.Range("E7").Value = "-4.10%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "8.206"
$ws.Range("E8").Value = "-1.38%"
$ws.Range("D9").Value = "2.933"
$ws.Range("E9").Value = "-10.90%"
$ws.Range("D10").Value = "0.9166"
$ws.Range("E10").Value = "-2.56%"
$ws.Range("D11").Value = "0.1208"
$ws.Range("E11").Value = "-8.26%"
$ws.Range("D12").Value = "0.1918"
$ws.Range("E12").Value = "-2.78%"
$ws.Range("D13").Value = "0.08928"
$ws.Range("E13").Value = "-0.84%"
$ws.Range("D14").Value = "0.03424"
$ws.Range("E14").Value = "-1.99%"
$ws.Range("D15").Value = "0.09704"
$ws.Range("D16").Value = "0.001374"
$ws.Range("E16").Value = "-2.99%"
$ws.Range("D17").Value = "0.005869"
$ws.Range("E17").Value = "-5.10%"
$ws.Range("D18").Value = "3.550"
$ws.Range("E18").Value = "-0.90%"
$ws.Range("E19").Value = "-1.78%"
$ws.Range("E20").Value = "0.33%"
$ws.Range("D21").Value = "0.1268"
$ws.Range("E21").Value = "-3.75%"
$ws.Range("D22").Value = "0.2589"
$ws.Range("E22").Value = "3.91%"
$ws.Range("D23").Value = "0.02103"
$ws.Range("E23").Value = "5,585.18%"
$ws.Range("D24").Value = "0.04367"
$ws.Range("E24").Value = "-0.11%"
$ws.Range("E25").Value = "-2.80%"
$ws.Range("D26").Value = "0.004247"
$ws.Range("E26").Value = "-10.30%"
$ws.Range("D27").Value = "0.0001299"
$ws.Range("E27").Value = "-66.68%"
$ws.Range("D39").Value = "0.02114"
$ws.Range("E39").Value = "-5.80%"
$ws.Range("D40").Value = "0.04953"
$ws.Range("E40").Value = "-5.74%"
$ws.Range("D41").Value = "0.007676"
$ws.Range("E41").Value = "0.86%"
$ws.Range("D42").Value = "0.009875"
$ws.Range("E42").Value = "-4.57%"
$ws.Range("D43").Value = "0.1342"
$ws.Range("E43").Value = "-3.73%"
$ws.Range("D44").Value = "0.002060"
$ws.Range("E44").Value = "-2.13%"
$ws.Range("D45").Value = "0.009600"
$ws.Range("E45").Value = "5.27%"
$ws.Range("D46").Value = "0.00006627"
$ws.Range("E46").Value = "-3.11%"
$ws.Range("E47").Value = "-0.29%"
$ws.Range("D48").Value = "0.003042"
$ws.Range("E48").Value = "0.78%"
$ws.Range("E50").Value = "-0.29%"
$ws.Range("E51").Value = "-0.29%"
